$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "** Fill out a little personal bio blurb"
$ws.Range("C13").Value = "** any other information you would like on your profile."

$ws.Range("D19").Select()
